$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Regenerate header labels: "..._old" -> "..._FV2410", "..._new" -> "..._FV2504" ---
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $v = $cell.Value()
    if ($v -ne $null) {
        $newValue = $v -replace '_old$', '_FV2410'
        $newValue = $newValue -replace '_new$', '_FV2504'
        if ($newValue -ne $v) {
            $cell.Value = $newValue
        }
    }
}

# --- Turn the sheet's used range into an Excel Table ("Table1") with a header row ---
$lastRow = $usedRange.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (split below row 1) ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
Write-Output "Done"
